$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column cells keep their original plain-text representation
# (avoids Excel auto-converting strings like "212.70" or "9.00" into numbers).

# Row 2: 'Bitcoin'
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.269.10"
$ws.Range("E2").Value = "  +0.30%  "

# Row 3: 'Ethereum'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.589.47"
$ws.Range("E3").Value = "  +0.52%  "

# Row 4: 'TetherUSD'
$ws.Range("E4").Value = "  -0.20%  "

# Row 5: 'BNB'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.70"
$ws.Range("E5").Value = "  +1.50%  "

# Row 6: 'XRP'
$ws.Range("E6").Value = "  +0.54%  "

# Row 7: 'USDC'
$ws.Range("E7").Value = "  -0.19%  "

# Row 8: 'Cardano'
$ws.Range("E8").Value = "  +0.13%  "

# Row 9: 'Dogecoin'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0608"
$ws.Range("E9").Value = "  -0.25%  "

# Row 10: 'Solana'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.33"
$ws.Range("E10").Value = "  -1.00%  "

# Row 11: 'TRON'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0849"
$ws.Range("E11").Value = "  +0.56%  "

# Row 12: 'WrappedliquidstakedEther2.0'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.813.31"
$ws.Range("E12").Value = "  +0.53%  "

# Row 13: 'WrappedEther'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.618.49"
$ws.Range("E13").Value = "  +2.24%  "

# Row 14: 'Polkadot'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.04"
$ws.Range("E14").Value = "  -0.18%  "

# Row 15: 'Polygon'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.520"
$ws.Range("E15").Value = "  +1.12%  "

# Row 16: 'Litecoin'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.40"
$ws.Range("E16").Value = "  -0.18%  "

# Row 17: 'WrappedBTC'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.283.88"
$ws.Range("E17").Value = "  +0.32%  "

# Row 18: 'ShibaInu'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0₃0726"
$ws.Range("E18").Value = "  -1.16%  "

# Row 19: 'Chainlink'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.45"
$ws.Range("E19").Value = "  +2.25%  "

# Row 20: 'BitcoinCash'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "213.11"
$ws.Range("E20").Value = "  +2.54%  "

# Row 21: 'Dai'
$ws.Range("E21").Value = "  -0.16%  "

# Row 22: 'Uniswap'
$ws.Range("E22").Value = "  +0.68%  "

# Row 23: 'Avalanche'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.00"
$ws.Range("E23").Value = "  +1.27%  "

# Row 24: 'Toncoin'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.16"
$ws.Range("E24").Value = "  -2.21%  "

# Row 25: 'Monero'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.69"
$ws.Range("E25").Value = "  +0.20%  "

# Row 26: 'BinanceUSD'
$ws.Range("E26").Value = "  -0.18%  "

# Row 27: 'Cosmos'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.04"
$ws.Range("E27").Value = "  +0.61%  "

# Row 28: 'Stellar'
$ws.Range("E28").Value = "  -0.53%  "

# Row 29: 'EthereumClassic'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.17"
$ws.Range("E29").Value = "  -0.29%  "

# Row 30: 'Hedera'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0499"
$ws.Range("E30").Value = "  -0.97%  "

# Row 31: 'PancakeSwap'
$ws.Range("E31").Value = "  +0.97%  "

# Row 32: 'Filecoin'
$ws.Range("E32").Value = "  -0.31%  "

# Row 33: 'InternetComputer(DFINITY)'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.95"
$ws.Range("E33").Value = "  +0.13%  "

# Row 34: 'Maker'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.337.37"
$ws.Range("E34").Value = "  +4.82%  "

# Row 35: 'HuobiToken'
$ws.Range("E35").Value = "  -0.87%  "

# Row 36: 'LidoDAOToken'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.47"
$ws.Range("E36").Value = "  -0.77%  "

# Row 37: 'ImmutableX'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.593"
$ws.Range("E37").Value = "  -2.84%  "

# Row 38: 'VeChain'
$ws.Range("E38").Value = "  -0.27%  "

# Row 39: 'ARBITRUM'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.819"
$ws.Range("E39").Value = "  +0.27%  "

# Row 40: 'FraxShare' -> 'WEMIXToken'
$ws.Range("B40").Value = "WEMIXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.03"
$ws.Range("E40").Value = "  -1.10%  "

# Row 41: 'PaxDollar' -> 'FraxShare'
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.73"
$ws.Range("E41").Value = "  +3.55%  "

# Row 42: 'WEMIXToken' -> 'PaxDollar'
$ws.Range("B42").Value = "PaxDollar"
$ws.Range("C42").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  -0.16%  "

# Row 43: 'MXToken'
$ws.Range("E43").Value = "  +0.26%  "

# Row 44: 'TrustWalletToken'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.763"
$ws.Range("E44").Value = "  -0.13%  "

# Row 45: 'Aave'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "61.84"
$ws.Range("E45").Value = "  -0.68%  "

# Row 46: 'RocketPoolETH'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.724.73"
$ws.Range("E46").Value = "  +0.42%  "

# Row 47: 'Quant'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "86.68"
$ws.Range("E47").Value = "  -2.67%  "

# Row 48: 'RenderToken'
$ws.Range("E48").Value = "  -3.36%  "

# Row 49: 'Algorand' -> 'Cronos'
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0503"
$ws.Range("E49").Value = "  -0.54%  "

# Row 50: 'Cronos' -> 'Algorand'
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0978"
$ws.Range("E50").Value = "  -2.45%  "

# Row 51: 'USDD'
$ws.Range("E51").Value = "  -0.48%  "
